# DemoFor.xlsx - 3rd commit
# Update the data row (row 2) with corrected/renamed values and append a
# new data row (row 3) that is a duplicate of row 2 but with an
# incremented "AddTag" value and the next calendar day.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Duplicate row 2 (values + formatting) into row 3 first, so the new
# row inherits the exact same cell styles (fonts / number formats) that
# row 2 currently has. ---
$ws.Range("A2:M2").Copy()
$ws.Range("A3:M3").PasteSpecial(-4122)   # xlPasteFormats

# Grab the existing "Notes" text (an ellipsis) so row 3 reuses the exact
# same shared string instead of creating a near-duplicate one.
$notesValue = $ws.Cells.Item(2, 12).Value2

# --- Correct / rename the text that is already in row 2 ---
$ws.Cells.Item(2, 1).Value2 = "Rakesh11"            # A2 SelectCompany
$ws.Cells.Item(2, 4).Value2 = "AutomationModel"     # D2 Model
$ws.Cells.Item(2, 6).Value2 = "AutomationAsset"     # F2 AssetName
$ws.Cells.Item(2, 8).Value2 = "AutomationSupplier"  # H2 Supplier
$ws.Cells.Item(2, 13).Value2 = "Parola"             # M2 DefaultLocation

# --- Fill in the new row 3 ---
$ws.Cells.Item(3, 1).Value2 = "Rakesh11"             # A3 SelectCompany
$ws.Cells.Item(3, 2).Value2 = "us-9877"              # B3 AddTag
$ws.Cells.Item(3, 3).Value2 = 7865                   # C3 Serial
$ws.Cells.Item(3, 4).Value2 = "AutomationModel"      # D3 Model
$ws.Cells.Item(3, 5).Value2 = "Ready to Deploy"      # E3 Status
$ws.Cells.Item(3, 6).Value2 = "AutomationAsset"      # F3 AssetName
$ws.Cells.Item(3, 7).Value2 = 44946                  # G3 CelenderDate
$ws.Cells.Item(3, 8).Value2 = "AutomationSupplier"   # H3 Supplier
$ws.Cells.Item(3, 9).Value2 = 1234                   # I3 OrderNo
$ws.Cells.Item(3, 10).Value2 = 35000                 # J3 PurchaseCost
$ws.Cells.Item(3, 11).Value2 = 6                     # K3 Warranty
$ws.Cells.Item(3, 12).Value2 = $notesValue           # L3 Notes
$ws.Cells.Item(3, 13).Value2 = "Parola"              # M3 DefaultLocation

# --- Update the sheet's active selection to reflect where editing ended ---
$ws.Activate() | Out-Null
$ws.Range("L4").Select() | Out-Null
